$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the survey's "label:English" header column to the xlsform-style
# "label::English (en)" column (matches the IETF language tag form already
# used on the choices sheet, so the shared-string table de-dupes them).
$ws.Range("C1").Value = "label::English (en)"

# Leaves the active cell on C2, mirroring the post-edit selection state.
[void]$ws.Range("C2").Select()
